# Updates loading_percent values for Case_2_69 (rows 2-25, sheet1) per commit "case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndexes = @(2, 3, 5, 6, 7, 8, 9, 11, 12, 13, 14, 15)

$data = @(
    @(10.33262446104749, 9.873561027317718, 13.00405383378493, 16.86991607391233, 28.20137921731904, 14.32643615764225, 22.66808073631774, 7.552965008185082, 10.05778039922445, 13.22454778644626, 19.09107358500509, 21.69902459128792),
    @(10.09266965951187, 9.876141332415109, 13.03341924948185, 15.89584955866808, 28.3245588910671, 14.37139436651309, 22.75696977219811, 7.375224935248127, 10.06439272228139, 13.18783595214496, 19.13897472483206, 21.78092948427687),
    @(9.943704471466837, 9.878151499105732, 13.05328376789502, 15.26997757108491, 28.407370458478, 14.40076975599295, 22.81497937680768, 7.264577680644237, 10.06978560556042, 13.16688245750329, 19.16993558302467, 21.83484237639931),
    @(9.88267535005374, 9.879078213649187, 13.0618402273675, 15.00819731993403, 28.4429174845152, 14.41318638712132, 22.83948252776553, 7.21916851446175, 10.07231912286326, 13.15874887215216, 19.182942934901, 21.85772350044326),
    @(9.872524482396516, 9.87923860372914, 13.06328890390565, 14.96433081551593, 28.44892866180529, 14.41527511026502, 22.84360345408443, 7.211610903360419, 10.07276011674933, 13.15742293148983, 19.18512641213715, 21.86157793193659),
    @(9.942882609715888, 9.878163560986817, 13.05339729425368, 15.26647399323137, 28.40784257339457, 14.40093540466928, 22.81530633627204, 7.263966487480182, 10.06981841263334, 13.16677111720774, 19.17010942212702, 21.83514726966354),
    @(10.25027707746877, 9.874362614908218, 13.01379857645203, 16.53996406344768, 28.24235893221665, 14.34157057885016, 22.69801804509341, 7.492030331981823, 10.05978413012583, 13.21156363738532, 19.10726890039798, 21.72651355033389),
    @(10.83655548192569, 9.870267442974004, 12.95068303718153, 19.00274580682531, 27.97500443315293, 14.23917915583404, 22.49520065803353, 8.033252703753064, 10.05064980230035, 13.3117286990523, 18.99628870347987, 21.54222227566351),
    @(11.25274438842421, 9.869279079210697, 12.91315361893745, 20.67494806633232, 27.81368994030746, 14.17246163711538, 22.36270475388878, 8.418933379144404, 10.05031878277117, 13.39245502167847, 18.9221579454228, 21.42433311739741),
    @(11.43806684933185, 9.869262348402483, 12.89799548058219, 21.3917225636224, 27.74798264632985, 14.14394969246939, 22.30600285414072, 8.587313226804779, 10.05154200688063, 13.43064210171907, 18.89002918414334, 21.37450250922416),
    @(11.50760344027402, 9.869317758266028, 12.89253030366823, 21.65686569030329, 27.72420912001829, 14.13341670743, 22.28504405482627, 8.650025345616282, 10.05220163818679, 13.44530533383967, 18.87809107747663, 21.35617904901484),
    @(11.4926569853929, 9.869303086665285, 12.89369510820745, 21.60004134736742, 27.72927980499056, 14.1356734457426, 22.28953510076072, 8.636566340229653, 10.05205085765852, 13.44213846222876, 18.88065201887973, 21.3601010334577),
    @(11.44380086899761, 9.869265672427908, 12.89754035042348, 21.4136618050453, 27.74600454375926, 14.14307785170242, 22.30426828326248, 8.592493803816462, 10.05159234529092, 13.43184444034222, 18.88904245690141, 21.37298407708377),
    @(11.4137897590786, 9.869250781221513, 12.89993145826397, 21.29868154950795, 27.7563934118627, 14.14764761432035, 22.31335957229667, 8.565360485808981, 10.0513370384026, 13.42556520102098, 18.89421155749079, 21.38094646503981),
    @(11.24054686830191, 9.869288840450011, 12.91418272572754, 20.62722412089977, 27.81813891966365, 14.17436190924619, 22.36648216957692, 8.40778419138846, 10.0502663920515, 13.38998819734031, 18.92428963490222, 21.42766610189482),
    @(11.13319444208167, 9.869422703853758, 12.92341542685701, 20.20408069597325, 27.85798703417264, 14.19122073737043, 22.39998542064902, 8.309280908587628, 10.04996067251793, 13.36853231713309, 18.94314916144975, 21.45729991796849),
    @(11.07107487507228, 9.869540473685198, 12.92890602865321, 19.95656407809801, 27.88162879704357, 14.20109053521898, 22.41959176907808, 8.25196084648652, 10.04991421256068, 13.35632977805968, 18.95414671597686, 21.47470195566292),
    @(11.04998021420429, 9.869587366918942, 12.93079601192033, 19.87204792380568, 27.88975738229833, 14.20446201417392, 22.42628789554281, 8.232440369346811, 10.04992073653978, 13.35222220196631, 18.95789609582379, 21.48065537337628),
    @(11.14466141383053, 9.869404237419669, 12.9224139430525, 20.24955283636154, 27.85367035698209, 14.18940818011266, 22.39638415467543, 8.319835683114139, 10.0499798342737, 13.37080207384161, 18.94112600859096, 21.45410835150804),
    @(11.45816898632342, 9.869274990169108, 12.89640345252965, 21.46857628470577, 27.74106196437187, 14.14089584044674, 22.29992687252276, 8.605467709080239, 10.05172170002478, 13.43486260974675, 18.88657179071522, 21.36918518655076),
    @(11.65930132154757, 9.86955018297661, 12.88100611446644, 22.22866616901552, 27.6739285375292, 14.11072805113351, 22.23987626722269, 8.786013790644239, 10.05400443207145, 13.47790700919878, 18.85224813224369, 21.31686740213822),
    @(11.55231778035409, 9.869370565349863, 12.88907751071895, 21.82633154458858, 27.70916606958522, 14.1266886054192, 22.27165301121129, 8.69022389418746, 10.05268178329994, 13.45482838152177, 18.87044582308692, 21.34449894525912),
    @(11.13947844062273, 9.86941245894546, 12.92286614508579, 20.22900810905287, 27.85561964527625, 14.19022708452609, 22.39801121270391, 8.31506601374466, 10.04997076840171, 13.36977550278912, 18.94204019324461, 21.45555012127429),
    @(10.68020074325863, 9.871018525173666, 12.96620318045849, 18.34778573295695, 28.04118387553025, 14.26538160921793, 22.54716398924598, 7.885199616085576, 10.05199676672479, 13.28334973983436, 19.0250067541705, 21.58900238452614)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $excelRow = $startRow + $i
    for ($j = 0; $j -lt $colIndexes.Length; $j++) {
        $ws.Cells.Item($excelRow, $colIndexes[$j]).Value = $rowValues[$j]
    }
}
